{"js": "// Insert a new \"Date\" styled paragraph containing \"2024-04-11\" right\n// after the Author paragraph (\"... Benjamin ABEL (facilitator)\").\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the paragraph with style \"Author\" (the one ending in \"(facilitator)\").\nlet authorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.style === \"Author\") {\n    authorParagraph = para;\n    break;\n  }\n}\n\nif (!authorParagraph) {\n  throw new Error('Could not find the \"Author\" paragraph to insert the date after.');\n}\n\nconst dateParagraph = authorParagraph.insertParagraph(\"2024-04-11\", Word.InsertLocation.after);\ndateParagraph.style = \"Date\";\n\nawait context.sync();\n", "ps1": "# Insert a new \"Date\" styled paragraph containing \"2024-04-11\" right\n# after the Author paragraph (\"... Benjamin ABEL (facilitator)\").\n\n$d = $word.ActiveDocument\n\n$authorParagraph = $null\nforeach ($para in $d.Paragraphs) {\n    if ($para.Style.NameLocal -eq \"Author\") {\n        $authorParagraph = $para\n        break\n    }\n}\n\nif ($authorParagraph -eq $null) {\n    throw \"Could not find the 'Author' paragraph to insert the date after.\"\n}\n\n$authorParagraph.Range.InsertParagraphAfter()\n\n$dateParagraph = $authorParagraph.Next()\n$dateParagraph.Range.Text = \"2024-04-11\"\n$dateParagraph.Style = \"Date\"\n"}
